# Insert a new weekly record row into the Acelga - Macroferia Regional de Talca sheet.
# This shifts the existing rows 315:349 down to 316:350, then fills the newly
# opened row 315 with the new observation's data (matching the pattern of the
# surrounding rows for the constant columns).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 315, pushing rows 315-349 down to 316-350.
$ws.Rows.Item(315).Insert()

# Fill in the new row 315 with the new weekly observation.
$ws.Range("A315").Value = 5
$ws.Range("B315").Value = "Macroferia Regional de Talca"
$ws.Range("C315").Value = "Maule"
$ws.Range("D315").Value = 44918
$ws.Range("E315").Value = 7
$ws.Range("F315").Value = 100112009
$ws.Range("G315").Value = "Acelga"
$ws.Range("H315").Value = "Sin especificar"
$ws.Range("I315").Value = "Primera"
$ws.Range("J315").Value = 500
$ws.Range("K315").Value = 2500
$ws.Range("L315").Value = 2500
$ws.Range("M315").Value = 2500
$ws.Range("N315").Value = "$/docena de atados (4 kilos)"
$ws.Range("O315").Value = "Región del Maule"
$ws.Range("P315").Value = 625
$ws.Range("Q315").Value = 4
$ws.Range("R315").Value = "Hortaliza"
